# Updates the cryptos list with fresh price/volume data.
# Forces cell values to be stored as text (matching the source data's
# inline-string representation) even when a value looks numeric
# (e.g. "189.23" or "1.00"), by temporarily applying a text number
# format, assigning the value, and then restoring a plain/default style
# so no stray formatting is left behind on the cell.
function Set-CellText {
    param($Worksheet, $CellRef, $Text)
    $cell = $Worksheet.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 'D2' '69.255.76'
Set-CellText $ws 'E2' '  +1.75%  '
Set-CellText $ws 'D3' '3.309.63'
Set-CellText $ws 'E3' '  +1.97%  '
Set-CellText $ws 'D5' '189.23'
Set-CellText $ws 'E5' '  +2.07%  '
Set-CellText $ws 'D6' '589.59'
Set-CellText $ws 'E6' '  +1.27%  '
Set-CellText $ws 'E7' '  -0.01%  '
Set-CellText $ws 'D8' '0.605'
Set-CellText $ws 'E8' '  +1.05%  '
Set-CellText $ws 'D9' '0.132'
Set-CellText $ws 'E9' '  +1.83%  '
Set-CellText $ws 'D10' '6.68'
Set-CellText $ws 'E10' '  +0.79%  '
Set-CellText $ws 'D11' '0.417'
Set-CellText $ws 'E11' '  -0.27%  '
Set-CellText $ws 'D12' '3.882.49'
Set-CellText $ws 'E12' '  +1.93%  '
Set-CellText $ws 'E13' '  +0.98%  '
Set-CellText $ws 'D14' '27.94'
Set-CellText $ws 'E14' '  -0.33%  '
Set-CellText $ws 'D15' '69.183.26'
Set-CellText $ws 'E15' '  +1.64%  '
Set-CellText $ws 'D16' '0.0000171'
Set-CellText $ws 'E16' '  +0.52%  '
Set-CellText $ws 'D17' '3.307.81'
Set-CellText $ws 'E17' '  +1.93%  '
Set-CellText $ws 'D18' '5.79'
Set-CellText $ws 'E18' '  -0.47%  '
Set-CellText $ws 'D19' '13.62'
Set-CellText $ws 'E19' '  +0.88%  '
Set-CellText $ws 'D20' '419.96'
Set-CellText $ws 'E20' '  +6.11%  '
Set-CellText $ws 'D21' '7.66'
Set-CellText $ws 'E21' '  +0.78%  '
Set-CellText $ws 'D22' '72.03'
Set-CellText $ws 'E22' '  +0.67%  '
Set-CellText $ws 'E23' '  +0.33%  '
Set-CellText $ws 'D24' '0.515'
Set-CellText $ws 'E24' '  -0.48%  '
Set-CellText $ws 'D25' '0.0000120'
Set-CellText $ws 'E25' '  +1.00%  '
Set-CellText $ws 'D26' '0.190'
Set-CellText $ws 'E26' '  +1.13%  '
Set-CellText $ws 'D27' '9.60'
Set-CellText $ws 'E27' '  -0.67%  '
Set-CellText $ws 'E28' '  +1.08%  '
Set-CellText $ws 'D29' '1.97'
Set-CellText $ws 'E29' '  -0.19%  '
Set-CellText $ws 'D30' '22.99'
Set-CellText $ws 'E30' '  +0.75%  '
Set-CellText $ws 'D31' '5.56'
Set-CellText $ws 'E31' '  -1.56%  '
Set-CellText $ws 'D32' '1.27'
Set-CellText $ws 'E32' '  +0.45%  '
Set-CellText $ws 'D33' '6.97'
Set-CellText $ws 'E33' '  -1.44%  '
Set-CellText $ws 'D34' '164.37'
Set-CellText $ws 'E34' '  +1.37%  '
Set-CellText $ws 'D35' '1.48'
Set-CellText $ws 'E35' '  -1.16%  '
Set-CellText $ws 'D36' '1.93'
Set-CellText $ws 'E36' '  +0.62%  '
Set-CellText $ws 'D37' '26.83'
Set-CellText $ws 'E37' '  +1.37%  '
Set-CellText $ws 'D38' '4.55'
Set-CellText $ws 'E38' '  -1.64%  '
Set-CellText $ws 'D39' '0.799'
Set-CellText $ws 'E39' '  -2.13%  '
Set-CellText $ws 'D40' '6.43'
Set-CellText $ws 'E40' '  -1.48%  '
Set-CellText $ws 'D41' '2.704.58'
Set-CellText $ws 'E41' '  +3.84%  '
Set-CellText $ws 'B42' 'dogwifhat'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText $ws 'D42' '2.47'
Set-CellText $ws 'E42' '  -0.95%  '
Set-CellText $ws 'B43' 'Hedera'
Set-CellText $ws 'C43' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText $ws 'D43' '0.0688'
Set-CellText $ws 'E43' '  +0.13%  '
Set-CellText $ws 'D44' '40.69'
Set-CellText $ws 'E44' '  -0.79%  '
Set-CellText $ws 'B45' 'InjectiveProtocol'
Set-CellText $ws 'C45' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-CellText $ws 'D45' '25.10'
Set-CellText $ws 'E45' '  -0.47%  '
Set-CellText $ws 'B46' 'Bittensor'
Set-CellText $ws 'C46' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-CellText $ws 'D46' '341.06'
Set-CellText $ws 'E46' '  +1.03%  '
Set-CellText $ws 'D47' '0.0279'
Set-CellText $ws 'E47' '  -0.26%  '
Set-CellText $ws 'D48' '32.31'
Set-CellText $ws 'E48' '  +3.50%  '
Set-CellText $ws 'D49' '1.00'
Set-CellText $ws 'E49' '  +2.11%  '
Set-CellText $ws 'D50' '6.26'
Set-CellText $ws 'E50' '  -1.46%  '
Set-CellText $ws 'E51' '  -0.73%  '
